$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.854.80"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.876.37"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.08%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7221"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.22"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.14%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3151"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.38%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07440"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.58"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.63%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08199"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.46%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7457"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.59%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.879.29"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.331"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.62"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.927.95"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.017"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "247.49"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.53%  "

$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.51"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.42%  "

$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007908"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.04%  "

$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.29%  "

$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.147.73"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.42%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.729"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.02%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1503"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.21"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.73%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.009"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.432"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.78%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.528"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.185"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.97%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05407"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.76%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.231"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7363"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.004"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.64%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.705"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.15%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.72%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4457"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.48%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8944"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.014"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.64"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.06%  "

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.043.11"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.63%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.002"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.78"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.45%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.467"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.34%  "

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.811"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.44%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.622"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.86%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.042.80"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.07%  "

